$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 35

$ws.Cells.Item($row, 1).Value = "Globo"
$ws.Cells.Item($row, 2).Value = "RJ TV 2"
$ws.Cells.Item($row, 3).Value = "Social"
$ws.Cells.Item($row, 4).Value = "2025-04-03T19:10"
$ws.Cells.Item($row, 5).Value = "Negativo"
$ws.Cells.Item($row, 6).Value = 'Segundo dia de mutirão do CadÚnico. Mais de 800 senhas foram distribuídas, mas alguns dos problemas continuam. Nos dois dias, reclamações são as mesmas: falta de estrutura, banheiros e água. Depoimentos de beneficiárias que estavam na fila reclamando e de uma moradora da rua, alegando "falta de respeito das autoridades campistas". Triagem começou de madrugada na fila. 800 vagas por dia. Depoimento de uma senhora alegando que conseguiu resolver seu problema. Novos depoimentos com reclamações. Entrevista com o assessor da secretaria, Ruan Barros. Dias 15 e 16, mutirão em Guarus; 29 e 30, mutirão na Baixada Campista. *matéria* Às 19h12, repórter *ao vivo* e atualização das informações com entrevista ao vido do assessor da secretaria. '

$wb.Save()
